$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" ---
# This shared string is used on the Overview sheet (E2, F2 - zh-cn/de-de
# status) and on each per-locale sheet's Status column (C2).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ("Ready for handoff" -eq $cell.Value2) {
            $cell.Value = "In Translation"
        }
    }
}

# --- Column width changes ---
# Overview sheet: columns E and F (zh-cn / de-de) narrower.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.577
$overview.Columns.Item(6).ColumnWidth = 12.577

# zh-cn sheet: column C (Status) narrower.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.577

# de-de sheet: column C (Status) narrower.
$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.577
